$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ASML)
$ws.Range("D2").Value = 1099.47
$ws.Range("E2").Value = 63.1
$ws.Range("F2").Value = 3.72
$ws.Range("K2").Value = 64.09999999999999
$ws.Range("N2").Value = 52.28493729186943

# Row 3 (TSM)
$ws.Range("D3").Value = 294.72
$ws.Range("E3").Value = 59.7
$ws.Range("F3").Value = 1.1
$ws.Range("H3").Value = 63
$ws.Range("K3").Value = 58.9
$ws.Range("N3").Value = 52.28493729186943

# Row 4 (AMD)
$ws.Range("D4").Value = 217.97
$ws.Range("E4").Value = 33.5
$ws.Range("F4").Value = 0.2
$ws.Range("K4").Value = 56.1
$ws.Range("N4").Value = 52.28493729186943

# Row 5 (NVDA)
$ws.Range("D5").Value = 182.41
$ws.Range("F5").Value = 3.06
$ws.Range("K5").Value = 51.1
$ws.Range("N5").Value = 52.28493729186943

# Row 6 (QCOM)
$ws.Range("D6").Value = 174.81
$ws.Range("E6").Value = 52.4
$ws.Range("F6").Value = 4.53
$ws.Range("I6").Value = 36
$ws.Range("K6").Value = 48.1
$ws.Range("N6").Value = 52.28493729186943
